$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This sheet stores every data cell as literal TEXT (t="inlineStr"), even
# values that look like integers/decimals/dates. Plain `Range.Value =
# "123"` assignment lets Excel's automatic type detection turn those into
# real numbers/dates, so instead we compute the literal text via a scratch
# formula cell and paste-special VALUES ONLY — that lands the result as a
# plain string with no numeric reinterpretation and without leaving any
# extra number-format/style baggage behind.
$scratch = $ws.Range("Z100")

function Set-TextCell($addr, $val) {
    $escaped = $val -replace '"', '""'
    $scratch.Formula = '="' + $escaped + '"'
    $scratch.Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4163) | Out-Null
}

# A truly empty (but present) inline-string cell can't be produced via the
# scratch-formula trick (pasting an empty value just clears the cell), so
# use a quote-prefixed empty entry instead, then strip the format it
# implicitly applies so the cell keeps the sheet's default (unstyled) look.
function Set-EmptyCell($addr) {
    $ws.Range($addr).Value = "'"
    $ws.Range($addr).Style = "Normal"
}

# --- New column L header ------------------------------------------------
$ws.Range("L1").Value = "Χαρακτηρισμός"
# Copy the header formatting (bold font + border + center/top alignment)
# from the neighbouring K1 header cell instead of re-declaring styles.
$ws.Range("K1").Copy() | Out-Null
$ws.Range("L1").PasteSpecial(-4122) | Out-Null

# --- Row 2: just gains the new (empty) L column --------------------------
Set-EmptyCell "L2"

# --- Row 3: several values change, plus the new L column -----------------
Set-TextCell "A3" "400008207899445"
Set-TextCell "F3" "2025-01-07"
Set-TextCell "I3" "70.03"
Set-TextCell "J3" "16.81"
Set-TextCell "K3" "86.84"
Set-EmptyCell "L3"

# --- Row 4: brand-new row --------------------------------------------------
Set-TextCell "A4" "400008195607600"
Set-EmptyCell "B4"
Set-EmptyCell "C4"
Set-TextCell "D4" "8Μ0ΤΔΑ"
Set-EmptyCell "E4"
Set-TextCell "F4" "2025-01-04"
Set-EmptyCell "G4"
Set-TextCell "H4" "1"
Set-TextCell "I4" "34.34"
Set-TextCell "J4" "8.24"
Set-TextCell "K4" "42.58"
Set-EmptyCell "L4"

# --- Row 5: brand-new row --------------------------------------------------
Set-TextCell "A5" "400008429648898"
Set-EmptyCell "B5"
Set-EmptyCell "C5"
Set-TextCell "D5" "8Μ0ΤΔΑ"
Set-EmptyCell "E5"
Set-TextCell "F5" "2025-01-29"
Set-EmptyCell "G5"
Set-TextCell "H5" "1"
Set-TextCell "I5" "146.07"
Set-TextCell "J5" "35.06"
Set-TextCell "K5" "181.13"
Set-EmptyCell "L5"

# Clean up the scratch cell used for the text-paste trick.
$scratch.ClearContents() | Out-Null
